$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 'disappointing'
$ws.Range("B4").Value = 0.8409090909090909
$ws.Range("C4").Value = 37
$ws.Range("D4").Value = 37
$ws.Range("H4").Value = 7
$ws.Range("K4").Value = 0.8923076923076924
$ws.Range("L4").Value = 58
$ws.Range("M4").Value = 58
$ws.Range("Q4").Value = 7
$ws.Range("A5").Value = 'poor'
$ws.Range("B5").Value = 0.7464788732394366
$ws.Range("C5").Value = 53
$ws.Range("D5").Value = 53
$ws.Range("H5").Value = 18
$ws.Range("K5").Value = 0.8172043010752689
$ws.Range("L5").Value = 76
$ws.Range("M5").Value = 76
$ws.Range("Q5").Value = 17
$ws.Range("A6").Value = 'disappointed'
$ws.Range("B6").Value = 0.7204301075268817
$ws.Range("C6").Value = 134
$ws.Range("D6").Value = 134
$ws.Range("H6").Value = 52
$ws.Range("K6").Value = 0.78125
$ws.Range("L6").Value = 50
$ws.Range("M6").Value = 50
$ws.Range("Q6").Value = 14
$ws.Range("A7").Value = 'however'
$ws.Range("B7").Value = 0.6875
$ws.Range("C7").Value = 44
$ws.Range("D7").Value = 44
$ws.Range("H7").Value = 20
$ws.Range("K7").Value = 0.7169811320754716
$ws.Range("L7").Value = 38
$ws.Range("M7").Value = 38
$ws.Range("Q7").Value = 15
$ws.Range("A8").Value = 'broke'
$ws.Range("B8").Value = 0.6407766990291263
$ws.Range("C8").Value = 132
$ws.Range("D8").Value = 132
$ws.Range("H8").Value = 74
$ws.Range("J8").Value = 'love'
$ws.Range("K8").Value = 0.5761494252873564
$ws.Range("L8").Value = 401
$ws.Range("M8").Value = 402
$ws.Range("P8").Value = $true
$ws.Range("Q8").Value = 295
$ws.Range("A9").Value = 'waste'
$ws.Range("B9").Value = 0.6351351351351351
$ws.Range("C9").Value = 94
$ws.Range("D9").Value = 94
$ws.Range("H9").Value = 54
$ws.Range("J9").Value = 'thank'
$ws.Range("K9").Value = 0.5652173913043478
$ws.Range("L9").Value = 39
$ws.Range("M9").Value = 39
$ws.Range("Q9").Value = 30
$ws.Range("A10").Value = 'instead'
$ws.Range("B10").Value = 0.625
$ws.Range("C10").Value = 30
$ws.Range("D10").Value = 30
$ws.Range("H10").Value = 18
$ws.Range("K10").Value = 0.5394190871369294
$ws.Range("L10").Value = 260
$ws.Range("M10").Value = 260
$ws.Range("Q10").Value = 222
$ws.Range("A11").Value = 'smaller'
$ws.Range("B11").Value = 0.5882352941176471
$ws.Range("C11").Value = 70
$ws.Range("D11").Value = 70
$ws.Range("H11").Value = 49
$ws.Range("K11").Value = 0.4675963904840033
$ws.Range("L11").Value = 570
$ws.Range("M11").Value = 571
$ws.Range("P11").Value = $true
$ws.Range("Q11").Value = 649
$ws.Range("A12").Value = 'guess'
$ws.Range("B12").Value = 0.5370370370370371
$ws.Range("C12").Value = 29
$ws.Range("D12").Value = 29
$ws.Range("H12").Value = 25
$ws.Range("K12").Value = 0.3577981651376147
$ws.Range("L12").Value = 117
$ws.Range("M12").Value = 117
$ws.Range("Q12").Value = 210
$ws.Range("A13").Value = 'small'
$ws.Range("B13").Value = 0.4898550724637681
$ws.Range("C13").Value = 169
$ws.Range("D13").Value = 169
$ws.Range("H13").Value = 176
$ws.Range("K13").Value = 0.3333333333333333
$ws.Range("L13").Value = 40
$ws.Range("M13").Value = 40
$ws.Range("Q13").Value = 80
$ws.Range("A14").Value = 'cheap'
$ws.Range("B14").Value = 0.4360189573459716
$ws.Range("C14").Value = 92
$ws.Range("D14").Value = 92
$ws.Range("H14").Value = 119
$ws.Range("K14").Value = 0.3313253012048193
$ws.Range("L14").Value = 55
$ws.Range("M14").Value = 55
$ws.Range("Q14").Value = 111
$ws.Range("B15").Value = 0.4330708661417323
$ws.Range("C15").Value = 55
$ws.Range("D15").Value = 55
$ws.Range("H15").Value = 72
$ws.Range("K15").Value = 0.3227513227513227
$ws.Range("L15").Value = 61
$ws.Range("M15").Value = 61
$ws.Range("Q15").Value = 128
$ws.Range("A16").Value = 'apart'
$ws.Range("B16").Value = 0.4315789473684211
$ws.Range("C16").Value = 41
$ws.Range("D16").Value = 41
$ws.Range("H16").Value = 54
$ws.Range("J16").Value = 'learn'
$ws.Range("K16").Value = 0.3046875
$ws.Range("L16").Value = 39
$ws.Range("M16").Value = 39
$ws.Range("Q16").Value = 89
$ws.Range("J17").Value = 'every'
$ws.Range("K17").Value = 0.232
$ws.Range("L17").Value = 29
$ws.Range("M17").Value = 29
$ws.Range("B18").Value = 0.4044943820224719
$ws.Range("C18").Value = 36
$ws.Range("D18").Value = 36
$ws.Range("H18").Value = 53
$ws.Range("J18").Value = 'christmas'
$ws.Range("K18").Value = 0.2208835341365462
$ws.Range("L18").Value = 55
$ws.Range("M18").Value = 55
$ws.Range("Q18").Value = 194
$ws.Range("A19").Value = 'ok'
$ws.Range("B19").Value = 0.3671875
$ws.Range("C19").Value = 47
$ws.Range("D19").Value = 47
$ws.Range("H19").Value = 81
$ws.Range("J19").Value = 'fun'
$ws.Range("K19").Value = 0.1980718667835232
$ws.Range("L19").Value = 226
$ws.Range("M19").Value = 226
$ws.Range("Q19").Value = 915
$ws.Range("B20").Value = 0.3469387755102041
$ws.Range("C20").Value = 34
$ws.Range("D20").Value = 34
$ws.Range("H20").Value = 64
$ws.Range("J20").Value = 'enjoy'
$ws.Range("K20").Value = 0.1935483870967742
$ws.Range("L20").Value = 36
$ws.Range("M20").Value = 36
$ws.Range("Q20").Value = 150
$ws.Range("B21").Value = 0.297029702970297
$ws.Range("C21").Value = 60
$ws.Range("D21").Value = 60
$ws.Range("H21").Value = 142
$ws.Range("K21").Value = 0.1116883116883117
$ws.Range("L21").Value = 172
$ws.Range("M21").Value = 173
$ws.Range("N21").Value = 0.99
$ws.Range("O21").Value = 0.01000000000000001
$ws.Range("Q21").Value = 1368
$ws.Range("A22").Value = 'size'
$ws.Range("B22").Value = 0.2422680412371134
$ws.Range("C22").Value = 47
$ws.Range("D22").Value = 47
$ws.Range("H22").Value = 147
$ws.Range("J22").Value = 'family'
$ws.Range("K22").Value = 0.1086350974930362
$ws.Range("L22").Value = 39
$ws.Range("M22").Value = 39
$ws.Range("Q22").Value = 320
$ws.Range("A23").Value = 'item'
$ws.Range("B23").Value = 0.2355072463768116
$ws.Range("C23").Value = 65
$ws.Range("D23").Value = 65
$ws.Range("H23").Value = 211
$ws.Range("J23").Value = 'easy'
$ws.Range("K23").Value = 0.09358288770053476
$ws.Range("L23").Value = 35
$ws.Range("M23").Value = 35
$ws.Range("Q23").Value = 339
$ws.Range("A24").Value = 'hard'
$ws.Range("B24").Value = 0.23
$ws.Range("C24").Value = 46
$ws.Range("D24").Value = 46
$ws.Range("H24").Value = 154
$ws.Range("K24").Value = 0.04806408544726302
$ws.Range("L24").Value = 36
$ws.Range("M24").Value = 39
$ws.Range("N24").Value = 0.92
$ws.Range("O24").Value = 0.07999999999999996
$ws.Range("Q24").Value = 713
$ws.Range("B25").Value = 0.2025316455696203
$ws.Range("C25").Value = 64
$ws.Range("D25").Value = 64
$ws.Range("H25").Value = 252
$ws.Range("A26").Value = 'price'
$ws.Range("B26").Value = 0.1936416184971098
$ws.Range("C26").Value = 67
$ws.Range("D26").Value = 69
$ws.Range("E26").Value = 0.03
$ws.Range("F26").Value = 0.97
$ws.Range("G26").Value = $true
$ws.Range("H26").Value = 279
$ws.Range("A27").Value = 'work'
$ws.Range("B27").Value = 0.1772151898734177
$ws.Range("C27").Value = 56
$ws.Range("D27").Value = 56
$ws.Range("H27").Value = 260
$ws.Range("A28").Value = 'box'
$ws.Range("B28").Value = 0.1727748691099476
$ws.Range("C28").Value = 33
$ws.Range("D28").Value = 33
$ws.Range("H28").Value = 158
$ws.Range("A29").Value = 'used'
$ws.Range("B29").Value = 0.1714285714285714
$ws.Range("C29").Value = 30
$ws.Range("D29").Value = 30
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = $false
$ws.Range("H29").Value = 145
$ws.Range("A30").Value = 'would'
$ws.Range("B30").Value = 0.1708766716196137
$ws.Range("C30").Value = 115
$ws.Range("D30").Value = 116
$ws.Range("E30").Value = 0.01
$ws.Range("F30").Value = 0.99
$ws.Range("H30").Value = 558
$ws.Range("A31").Value = 'better'
$ws.Range("B31").Value = 0.1588785046728972
$ws.Range("C31").Value = 34
$ws.Range("D31").Value = 34
$ws.Range("H31").Value = 180
$ws.Range("A32").Value = 'product'
$ws.Range("B32").Value = 0.1541850220264317
$ws.Range("C32").Value = 70
$ws.Range("D32").Value = 70
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = $false
$ws.Range("H32").Value = 384
$ws.Range("A33").Value = '3'
$ws.Range("B33").Value = 0.1330645161290323
$ws.Range("C33").Value = 33
$ws.Range("D33").Value = 33
$ws.Range("H33").Value = 215
$ws.Range("A34").Value = '2'
$ws.Range("B34").Value = 0.1283018867924528
$ws.Range("C34").Value = 34
$ws.Range("D34").Value = 36
$ws.Range("E34").Value = 0.06
$ws.Range("F34").Value = 0.9399999999999999
$ws.Range("G34").Value = $true
$ws.Range("H34").Value = 231
$ws.Range("A35").Value = 'little'
$ws.Range("B35").Value = 0.1146067415730337
$ws.Range("C35").Value = 51
$ws.Range("D35").Value = 55
$ws.Range("E35").Value = 0.07000000000000001
$ws.Range("F35").Value = 0.9299999999999999
$ws.Range("H35").Value = 394
$ws.Range("A36").Value = 'use'
$ws.Range("B36").Value = 0.1129476584022039
$ws.Range("C36").Value = 41
$ws.Range("D36").Value = 43
$ws.Range("E36").Value = 0.05
$ws.Range("F36").Value = 0.95
$ws.Range("H36").Value = 322
$ws.Range("B37").Value = 0.09014084507042254
$ws.Range("C37").Value = 32
$ws.Range("D37").Value = 32
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = $false
$ws.Range("H37").Value = 323
$ws.Range("B38").Value = 0.06754530477759473
$ws.Range("C38").Value = 41
$ws.Range("D38").Value = 42
$ws.Range("E38").Value = 0.02
$ws.Range("F38").Value = 0.98
$ws.Range("H38").Value = 566
$ws.Range("A39").Value = 'toy'
$ws.Range("B39").Value = 0.06605222734254992
$ws.Range("C39").Value = 43
$ws.Range("D39").Value = 47
$ws.Range("E39").Value = 0.09
$ws.Range("F39").Value = 0.91
$ws.Range("H39").Value = 608
$ws.Range("A40").Value = 'one'
$ws.Range("B40").Value = 0.05761843790012804
$ws.Range("C40").Value = 45
$ws.Range("D40").Value = 58
$ws.Range("E40").Value = 0.22
$ws.Range("F40").Value = 0.78
$ws.Range("H40").Value = 736

$ws.Rows.Item(41).Delete()

